# Apply the recorded changes to the "Ref_coal phase out" QC formula sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3: Year values M3:W3 changed from 2021 -> 2012 ---
$ws.Range("M3:W3").Value = 2012

# --- Row 4: Retrofit_existing_min schedule updated; trailing formulas
#     (=S4, =T4, =U4, =V4) replaced by plain literal values of 1 ---
$ws.Range("O4").Value = 0.2
$ws.Range("P4").Value = 1
$ws.Range("Q4").Value = 1
$ws.Range("R4").Value = 1
$ws.Range("S4").Value = 1
$ws.Range("T4").Value = 1
$ws.Range("U4").Value = 1
$ws.Range("V4").Value = 1
$ws.Range("W4").Value = 1

# --- Row 5: Year values M5:W5 changed from 2021 -> 2012 ---
$ws.Range("M5:W5").Value = 2012

# --- Row 6: Retrofit_existing_min schedule updated; trailing formulas
#     (=S6, =T6, =U6, =V6) replaced by plain literal values of 1 ---
$ws.Range("O6").Value = 0.2
$ws.Range("P6").Value = 1
$ws.Range("Q6").Value = 1
$ws.Range("R6").Value = 1
$ws.Range("S6").Value = 1
$ws.Range("T6").Value = 1
$ws.Range("U6").Value = 1
$ws.Range("V6").Value = 1
$ws.Range("W6").Value = 1

# --- Update the visible selection to match the author's final selection ---
$ws.Range("M3:W6").Select()
